$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.140.93"
$ws.Range("E2").Value = "  +0.99%  "

$ws.Range("D3").Value = "'1.639.36"
$ws.Range("E3").Value = "  -0.03%  "

$ws.Range("D5").Value = "'216.67"

$ws.Range("E6").Value = "  +2.31%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").Value = "'0.254"
$ws.Range("E8").Value = "  -0.22%  "

$ws.Range("D9").Value = "'0.0625"
$ws.Range("E9").Value = "  +0.03%  "

$ws.Range("D10").Value = "'19.95"
$ws.Range("E10").Value = "  +0.38%  "

$ws.Range("D11").Value = "'0.0847"
$ws.Range("E11").Value = "  -0.16%  "

$ws.Range("D12").Value = "'1.868.61"
$ws.Range("E12").Value = "  -0.03%  "

$ws.Range("D13").Value = "'1.630.82"
$ws.Range("E13").Value = "  -0.47%  "

$ws.Range("E14").Value = "  +0.19%  "

$ws.Range("E15").Value = "  +2.09%  "

$ws.Range("D16").Value = "'66.77"
$ws.Range("E16").Value = "  -0.78%  "

$ws.Range("D17").Value = "'27.144.95"
$ws.Range("E17").Value = "  +1.03%  "

$ws.Range("E18").Value = "  +1.25%  "

$ws.Range("D19").Value = "'217.07"
$ws.Range("E19").Value = "  -1.23%  "

$ws.Range("E22").Value = "  +3.50%  "

$ws.Range("E23").Value = "  +0.55%  "

$ws.Range("E24").Value = "  -0.31%  "

$ws.Range("D25").Value = "'147.11"
$ws.Range("E25").Value = "  -0.04%  "

$ws.Range("E26").Value = "  -0.13%  "

$ws.Range("D27").Value = "'7.42"
$ws.Range("E27").Value = "  +1.06%  "

$ws.Range("E28").Value = "  +0.00%  "

$ws.Range("E29").Value = "  -0.78%  "

$ws.Range("D30").Value = "'0.0509"
$ws.Range("E30").Value = "  +1.08%  "

$ws.Range("E31").Value = "  -0.16%  "

$ws.Range("E32").Value = "  +1.37%  "

$ws.Range("E33").Value = "  +0.53%  "

$ws.Range("D34").Value = "'1.307.52"
$ws.Range("E34").Value = "  +3.32%  "

$ws.Range("E35").Value = "  -0.03%  "

$ws.Range("E36").Value = "  +1.25%  "

$ws.Range("E37").Value = "  -1.41%  "

$ws.Range("D38").Value = "'0.857"
$ws.Range("E38").Value = "  +2.72%  "

$ws.Range("D39").Value = "'0.543"
$ws.Range("E39").Value = "  +1.76%  "

$ws.Range("E40").Value = "  +0.01%  "

$ws.Range("D41").Value = "'0.810"
$ws.Range("E41").Value = "  -0.16%  "

$ws.Range("D43").Value = "'5.30"
$ws.Range("E43").Value = "  -1.69%  "

$ws.Range("D44").Value = "'1.778.93"
$ws.Range("E44").Value = "  -0.03%  "

$ws.Range("D45").Value = "'61.77"
$ws.Range("E45").Value = "  -0.17%  "

$ws.Range("D46").Value = "'91.54"
$ws.Range("E46").Value = "  -0.29%  "

$ws.Range("E47").Value = "  +0.84%  "

$ws.Range("E48").Value = "  +1.85%  "

$ws.Range("E49").Value = "  -0.21%  "

$ws.Range("D50").Value = "'7.65"
$ws.Range("E50").Value = "  +0.07%  "

$ws.Range("D51").Value = "'0.0962"
$ws.Range("E51").Value = "  +0.02%  "
